$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.2219512195121951
$ws.Range("C2").Value = 0.5048780487804878
$ws.Range("J2").Value = 0.01951219512195122
$ws.Range("P2").Value = 0.1512195121951219
$ws.Range("S2").Value = 0.1024390243902439
$ws.Range("B3").Value = 0.01869158878504673
$ws.Range("C3").Value = 0.02336448598130841
$ws.Range("J3").Value = 0.04205607476635514
$ws.Range("P3").Value = 0.7616822429906542
$ws.Range("S3").Value = 0.1542056074766355
$ws.Range("J4").Value = 0.03076923076923077
$ws.Range("P4").Value = 0.6615384615384615
$ws.Range("S4").Value = 0.3076923076923077
$ws.Range("B6").Value = 0.0830188679245283
$ws.Range("D6").Value = 0.02264150943396226
$ws.Range("F6").Value = 0.1018867924528302
$ws.Range("J6").Value = 0.2037735849056604
$ws.Range("O6").Value = 0.02264150943396226
$ws.Range("Q6").Value = 0.2113207547169811
$ws.Range("R6").Value = 0.03018867924528302
$ws.Range("S6").Value = 0.3245283018867924
$ws.Range("B7").Value = 0.1197604790419162
$ws.Range("D7").Value = 0.01796407185628742
$ws.Range("F7").Value = 0.04191616766467066
$ws.Range("J7").Value = 0.1317365269461078
$ws.Range("O7").Value = 0.02395209580838323
$ws.Range("Q7").Value = 0.2275449101796407
$ws.Range("R7").Value = 0.08982035928143713
$ws.Range("S7").Value = 0.3473053892215569
$ws.Range("B8").Value = 0.1143984220907298
$ws.Range("D8").Value = 0.03155818540433925
$ws.Range("F8").Value = 0.07692307692307693
$ws.Range("J8").Value = 0.09270216962524655
$ws.Range("O8").Value = 0.02761341222879684
$ws.Range("Q8").Value = 0.1932938856015779
$ws.Range("R8").Value = 0.0670611439842209
$ws.Range("S8").Value = 0.3964497041420119
$ws.Range("B9").Value = 0.1371428571428571
$ws.Range("D9").Value = 0.03428571428571429
$ws.Range("F9").Value = 0.05714285714285714
$ws.Range("J9").Value = 0.08571428571428572
$ws.Range("O9").Value = 0.05142857142857143
$ws.Range("Q9").Value = 0.1942857142857143
$ws.Range("R9").Value = 0.05714285714285714
$ws.Range("S9").Value = 0.3828571428571428
$ws.Range("B10").Value = 0.13196894848271
$ws.Range("D10").Value = 0.02399435426958363
$ws.Range("E10").Value = 0.0007057163020465773
$ws.Range("F10").Value = 0.06422018348623854
$ws.Range("J10").Value = 0.1023288637967537
$ws.Range("O10").Value = 0.03952011291460833
$ws.Range("Q10").Value = 0.2342978122794637
$ws.Range("R10").Value = 0.05998588567395907
$ws.Range("S10").Value = 0.3429781227946366
$ws.Range("G11").Value = 0.1126760563380282
$ws.Range("J11").Value = 0.1091549295774648
$ws.Range("K11").Value = 0.1866197183098592
$ws.Range("L11").Value = 0.5598591549295775
$ws.Range("S11").Value = 0.03169014084507042
$ws.Range("G12").Value = 0.725
$ws.Range("J12").Value = 0.24375
$ws.Range("K12").Value = 0.00625
$ws.Range("L12").Value = 0.0125
$ws.Range("S12").Value = 0.0125
$ws.Range("G13").Value = 0.6666666666666666
$ws.Range("J13").Value = 0.303030303030303
$ws.Range("S13").Value = 0.0303030303030303
$ws.Range("F15").Value = 0.02912621359223301
$ws.Range("H15").Value = 0.1715210355987055
$ws.Range("I15").Value = 0.03559870550161812
$ws.Range("J15").Value = 0.3203883495145631
$ws.Range("K15").Value = 0.09061488673139159
$ws.Range("M15").Value = 0.003236245954692557
$ws.Range("N15").Value = 0.003236245954692557
$ws.Range("O15").Value = 0.08414239482200647
$ws.Range("S15").Value = 0.2621359223300971
$ws.Range("F16").Value = 0.03435114503816794
$ws.Range("H16").Value = 0.1908396946564886
$ws.Range("I16").Value = 0.08015267175572519
$ws.Range("J16").Value = 0.4389312977099237
$ws.Range("K16").Value = 0.08396946564885496
$ws.Range("M16").Value = 0.01526717557251908
$ws.Range("O16").Value = 0.0648854961832061
$ws.Range("S16").Value = 0.0916030534351145
$ws.Range("F17").Value = 0.03963963963963964
$ws.Range("H17").Value = 0.1675675675675676
$ws.Range("I17").Value = 0.07747747747747748
$ws.Range("J17").Value = 0.4738738738738739
$ws.Range("K17").Value = 0.07747747747747748
$ws.Range("M17").Value = 0.01261261261261261
$ws.Range("O17").Value = 0.06486486486486487
$ws.Range("S17").Value = 0.08648648648648649
$ws.Range("F18").Value = 0.03973509933774835
$ws.Range("H18").Value = 0.1788079470198675
$ws.Range("I18").Value = 0.1059602649006623
$ws.Range("J18").Value = 0.3774834437086093
$ws.Range("K18").Value = 0.07947019867549669
$ws.Range("M18").Value = 0.02649006622516556
$ws.Range("O18").Value = 0.0728476821192053
$ws.Range("S18").Value = 0.119205298013245
$ws.Range("F19").Value = 0.01320901320901321
$ws.Range("H19").Value = 0.2253302253302253
$ws.Range("I19").Value = 0.06682206682206682
$ws.Range("J19").Value = 0.3986013986013986
$ws.Range("K19").Value = 0.09557109557109557
$ws.Range("M19").Value = 0.01476301476301476
$ws.Range("N19").Value = 0.001554001554001554
$ws.Range("O19").Value = 0.08003108003108003
$ws.Range("S19").Value = 0.1041181041181041

Write-Output "Applied 109 cell updates"
